## Auto commit - 02041140
## Applies the 202602_HL_Maintain_Report.xlsx edit:
##  - bumps the report's "製表日期" (report-generation date) in the title cell A1
##    from 2026-02-03 to 2026-02-04
##  - extends the Print_Area named range from row 16 to row 17
##  - switches P16 / AC16 (work-content free-text cells) to a wrapping style
##  - appends a new data row (row 17, item #15) describing a new maintenance ticket
##  - moves the active selection to the newly added row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Print area: $A$1:$AK$16  ->  $A$1:$AK$17
# ---------------------------------------------------------------------------
$ws.PageSetup.PrintArea = "'Report'!`$A`$1:`$AK`$17"

# ---------------------------------------------------------------------------
# 2. Title cell: report date 2026-02-03 -> 2026-02-04
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "萊爾富 工作統計表  篩選月份：202602   (  製表日期:2026-02-04  )"

# ---------------------------------------------------------------------------
# 3. P16 / AC16 pick up the wrap-text variant of their existing style
# ---------------------------------------------------------------------------
$ws.Range("P16").WrapText = $true
$ws.Range("AC16").WrapText = $true

# ---------------------------------------------------------------------------
# 4. New row 17 - start by cloning row 15's formatting pattern (A:AK) so the
#    new row's cell styles line up with the surrounding striped rows without
#    spilling formatting past column AK.
# ---------------------------------------------------------------------------
$ws.Range("A15:AK15").Copy()
$ws.Range("A17:AK17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5. Populate the new row's values
# ---------------------------------------------------------------------------
$ws.Range("A17").Value  = 15
$ws.Range("B17").Value  = "維修"
$ws.Range("C17").Value  = 2026020765
$ws.Range("D17").Value  = "E3052115020401"
$ws.Range("E17").Value  = "一般件"
$ws.Range("F17").Value  = 3052
$ws.Range("G17").Value  = "新莊莊玲店"
$ws.Range("H17").Value  = "新北市新莊區"
$ws.Range("I17").Value  = "2026-02-04 09:48:39"
$ws.Range("J17").Value  = "星期三"
$ws.Range("K17").Value  = "上午"
$ws.Range("L17").Value  = "HLF3"
$ws.Range("M17").Value  = "HL-LIFE-ET QRcode掃描器"
$ws.Range("N17").Value  = "F301"
$ws.Range("O17").Value  = "掃描無反應或感應不良"
$ws.Range("P17").Value  = "門市反應MMK 四代機 QRCODE掃描器刷讀QRCODE有亮紅光但無反應(例:餐食券...)，已有重新開機仍異常(掃描後無反應也沒出紙)"
$ws.Range("Q17").Value  = "THILF03052"
$ws.Range("R17").Value  = "新北一"
$ws.Range("S17").Value  = "湯家瑋"
$ws.Range("T17").Value  = 1
$ws.Range("U17").Value  = "已完工"
$ws.Range("V17").Value  = "2026-02-04 09:58:52"
$ws.Range("W17").Value  = "2026-02-04 10:30:00"
$ws.Range("X17").Value  = "2026-02-04 11:00:00"
$ws.Range("Y17").Value  = "2026-02-05 13:58:00"
$ws.Range("Z17").Value  = 0.5
$ws.Range("AB17").Value = "到場處理"
$ws.Range("AC17").Value = "現場測試正常"
$ws.Range("AK17").Value = "O"

# ---------------------------------------------------------------------------
# 6. Move the selection to the new row, like the source workbook shows
# ---------------------------------------------------------------------------
$ws.Range("A17").Select()
